$d = $word.ActiveDocument

# Map of old exact text -> new exact text, applied as whole-string literal
# Find/Replace (MatchCase=true, MatchWholeWord=false, Wrap=FindContinue).
$replacements = @(
    @("2023-11-15 Wednesday", "2023-11-16 Thursday"),
    @("76×28=2128", "13×98=1274"),
    @("82×62=5084", "47×24=1128"),
    @("66×39=2574", "71×60=4260"),
    @("64×82=5248", "30×72=2160"),
    @("15×18=270", "15×31=465"),
    @("17×38=646", "20×24=480"),
    @("94×38=3572", "65×68=4420"),
    @("31×28=868", "56×48=2688"),
    @("16×77=1232", "82×79=6478"),
    @("73×19=1387", "35×92=3220"),
    @("40×28=1120", "95×80=7600"),
    @("40×68=2720", "81×34=2754"),
    @("58×60=3480", "40×79=3160"),
    @("32×52=1664", "49×31=1519"),
    @("36×14=504", "76×14=1064"),
    @("19×67=1273", "12×14=168"),
    @("62×80=4960", "14×27=378"),
    @("16×94=1504", "74×38=2812"),
    @("61×68=4148", "37×63=2331"),
    @("15×76=1140", "71×52=3692"),
    @("39×98=3822", "91×40=3640"),
    @("63×29=1827", "27×97=2619"),
    @("81×92=7452", "95×82=7790"),
    @("92×23=2116", "74×26=1924"),
    @("62×55=3410", "70×34=2380")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "WARNING: not found -> $old"
    }
}

Write-Output "Replacements applied"
